$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values to reflect repulled data / recalculated mean
$ws.Range("F6").Value = 3
$ws.Range("F8").Value = 6
$ws.Range("F10").Value = -3
$ws.Range("F16").Value = 3
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = -3
$ws.Range("F25").Value = -3
$ws.Range("F30").Value = -4
$ws.Range("F33").Value = -11
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = -3
$ws.Range("F55").Value = 1
$ws.Range("F56").Value = 5
$ws.Range("F57").Value = -14
$ws.Range("F59").Value = -1
